# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" (i.e. right before the
#    existing "2021-Q2" sheet), and fill it with the fund-holding table.
# 2) On the "总计" summary sheet, overwrite row 2 with the new 2022-Q4 totals
#    and append a new row 3 that carries the data which used to live in row 2
#    (the 2021-Q2 totals), re-numbering the index column.
#
# Note: this engine's `.Value` *getter* is unreliable (returns a placeholder
# string instead of the cached value), so reads use `.Value2` / `.Formula`
# instead; `.Value` is fine as a *setter* and is used for writes.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# --- 1. New "2022-Q4" worksheet, inserted before "2021-Q2" -----------------

$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q4"

# Header row (B1:H1)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Apply the bold/bordered header style (same xf the "总计" header row uses)
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# B2:G8 hold fund code / name / size / position values that are stored as
# plain text even though several look numeric ("12.13", "005299", ...).
# Mark the range as Text up front so the values keep leading zeros / exact
# digits instead of being auto-coerced to numbers.
$newSheet.Range("B2:G8").NumberFormat = "@"

$fundData = @(
  @(0, "005299", "万家成长优选灵活配置混合A", "12.13", "91.35", "3.62", "0.4391", 7),
  @(1, "005300", "万家成长优选灵活配置混合C", "9.48", "91.35", "3.62", "0.3432", 7),
  @(2, "010694", "万家内需增长一年持有期混合", "9.46", "94.46", "3.11", "0.2942", 8),
  @(3, "006132", "万家智造优势混合A", "4.10", "93.92", "4.08", "0.1673", 5),
  @(4, "006133", "万家智造优势混合C", "0.78", "93.92", "4.08", "0.0318", 5),
  @(5, "014831", "兴银中证1000指数增强A", "1.40", "82.60", "0.92", "0.0129", 7),
  @(6, "014832", "兴银中证1000指数增强C", "1.01", "82.60", "0.92", "0.0093", 7)
)

$row = 2
foreach ($rec in $fundData) {
  $newSheet.Cells.Item($row, 1).Value = $rec[0]
  $newSheet.Cells.Item($row, 2).Value = $rec[1]
  $newSheet.Cells.Item($row, 3).Value = $rec[2]
  $newSheet.Cells.Item($row, 4).Value = $rec[3]
  $newSheet.Cells.Item($row, 5).Value = $rec[4]
  $newSheet.Cells.Item($row, 6).Value = $rec[5]
  $newSheet.Cells.Item($row, 7).Value = $rec[6]
  $newSheet.Cells.Item($row, 8).Value = $rec[7]
  $row = $row + 1
}

# A2:A8 get the same index-column style as "总计"!A2 (bold/bordered, centered)
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

# --- 2. Update the "总计" summary sheet -------------------------------------

# Preserve the old row-2 (2021-Q2) values before they get overwritten. The
# `.Value` getter is unreliable here, so read back through `.Formula`.
$oldLabel = $totalSheet.Range("B2").Formula
$oldCount = $totalSheet.Range("C2").Formula
$oldValue = $totalSheet.Range("D2").Formula

# New row 2: 2022-Q4 totals
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 1.3

# New row 3: the previous 2021-Q2 totals, with the index column re-numbered
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Formula = $oldLabel
$totalSheet.Range("C3").Formula = $oldCount
$totalSheet.Range("D3").Formula = $oldValue

# A3 gets the same style as A2 (bold/bordered, centered)
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
